# Timesheet update: add a new punch-in entry for 2023-12-20 (row 21)
# covering 5 minutes of LO1 (utilities) work and 15 minutes of LO3
# (VFX+SFX+animation) work - the "blend animation" punch fix - and let the
# existing fill-down formulas (F column totals, placeholder (1/60)*(0)
# entries) extend one row further to row 22, matching Excel's behaviour
# when a value is entered in the row right below a filled table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New timesheet row: date + the two logged durations.
$ws.Range("A21").Value = 45280
$ws.Range("B21").Formula = "=(1/60)*(5)"
$ws.Range("D21").Formula = "=(1/60)*(15)"

# Excel auto-extends the adjacent "fill" formulas down into row 22 once
# row 21 stops being the blank trailing row - recreate that here.
$ws.Range("B22:E22").Formula = "=(1/60)*(0)"
$ws.Range("F22").Formula = "=SUM(B22:E22)"

# Leave the selection where the user ended up after entering the data.
$ws.Range("F22").Select()
